$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 3 that duplicates the current row 2 (the original PNG entry)
$ws.Rows.Item(3).Insert()

$ws.Range("A3").Value = $ws.Range("A2").Value()
$ws.Range("B3").Value = "C:/Users/Asus/Desktop/12 — копия\51411030m\007_image_51411030m.png"
$ws.Range("C3").Value = $ws.Range("C2").Value()
$ws.Range("D3").Value = $ws.Range("D2").Value()
$ws.Range("E3").Value = $ws.Range("E2").Value()

# Update the original row 2 to point to the jpg variant
$ws.Range("B2").Value = "C:/Users/Asus/Desktop/12 — копия\51411030m\007_image_51411030m.jpg"
